# Auto-generated PowerShell Excel COM-interop edit script
# Applies the betexplorer saudi-professional-league 2023-2024 update:
#  - rows 22-24, 83-84, 87-88, 92-93, 95-96: re-sequenced match rows
#    (same A:E index/date, F:V match data re-ordered between rows)
#  - rows 104-106: three newly scraped matches appended at the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Re-sequenced rows: only columns F:V (match data) change;
# columns A:E (Indice/pais/torneio/temporada/data_partida) stay put -----
# Row 22
$ws.Cells.Item(22,6).Value = 'Al Raed'
$ws.Cells.Item(22,7).Value = 0
$ws.Cells.Item(22,8).Value = 'Al Hilal'
$ws.Cells.Item(22,9).Value = 4
$ws.Cells.Item(22,10).Value = 12.64
$ws.Cells.Item(22,11).Value = '22/08/2023 07:46'
$ws.Cells.Item(22,12).Value = 18.25
$ws.Cells.Item(22,13).Value = '24/08/2023 19:59'
$ws.Cells.Item(22,14).Value = 6.75
$ws.Cells.Item(22,15).Value = '22/08/2023 07:46'
$ws.Cells.Item(22,16).Value = 9.07
$ws.Cells.Item(22,17).Value = '24/08/2023 19:59'
$ws.Cells.Item(22,18).Value = 1.21
$ws.Cells.Item(22,19).Value = '22/08/2023 07:46'
$ws.Cells.Item(22,20).Value = 1.14
$ws.Cells.Item(22,21).Value = '24/08/2023 19:51'
$ws.Cells.Item(22,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-raed-al-hilal/MN4PHx3L/'

# Row 23
$ws.Cells.Item(23,6).Value = 'Al Ettifaq'
$ws.Cells.Item(23,7).Value = 1
$ws.Cells.Item(23,8).Value = 'Al Khaleej'
$ws.Cells.Item(23,9).Value = 1
$ws.Cells.Item(23,10).Value = 1.56
$ws.Cells.Item(23,11).Value = '22/08/2023 07:46'
$ws.Cells.Item(23,12).Value = 1.85
$ws.Cells.Item(23,13).Value = '24/08/2023 19:54'
$ws.Cells.Item(23,14).Value = 4.25
$ws.Cells.Item(23,15).Value = '22/08/2023 07:46'
$ws.Cells.Item(23,16).Value = 3.7
$ws.Cells.Item(23,17).Value = '24/08/2023 19:54'
$ws.Cells.Item(23,18).Value = 4.92
$ws.Cells.Item(23,19).Value = '22/08/2023 07:46'
$ws.Cells.Item(23,20).Value = 4.2
$ws.Cells.Item(23,21).Value = '24/08/2023 19:54'
$ws.Cells.Item(23,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ettifaq-fc-al-khaleej/Gp4TGdIR/'

# Row 24
$ws.Cells.Item(24,6).Value = 'Al Riyadh'
$ws.Cells.Item(24,7).Value = 0
$ws.Cells.Item(24,8).Value = 'Al Ittihad'
$ws.Cells.Item(24,9).Value = 4
$ws.Cells.Item(24,10).Value = 13.51
$ws.Cells.Item(24,11).Value = '22/08/2023 07:46'
$ws.Cells.Item(24,12).Value = 12.1
$ws.Cells.Item(24,13).Value = '24/08/2023 19:54'
$ws.Cells.Item(24,14).Value = 6.52
$ws.Cells.Item(24,15).Value = '22/08/2023 07:46'
$ws.Cells.Item(24,16).Value = 6.94
$ws.Cells.Item(24,17).Value = '24/08/2023 19:54'
$ws.Cells.Item(24,18).Value = 1.21
$ws.Cells.Item(24,19).Value = '22/08/2023 07:46'
$ws.Cells.Item(24,20).Value = 1.22
$ws.Cells.Item(24,21).Value = '24/08/2023 19:07'
$ws.Cells.Item(24,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-riyadh-al-ittihad/zVp0Bztk/'

# Row 83
$ws.Cells.Item(83,6).Value = 'Al Akhdoud'
$ws.Cells.Item(83,7).Value = 1
$ws.Cells.Item(83,8).Value = 'Al Feiha'
$ws.Cells.Item(83,9).Value = 2
$ws.Cells.Item(83,10).Value = 3.22
$ws.Cells.Item(83,11).Value = '15/10/2023 17:20'
$ws.Cells.Item(83,12).Value = 2.86
$ws.Cells.Item(83,13).Value = '20/10/2023 16:58'
$ws.Cells.Item(83,14).Value = 3.32
$ws.Cells.Item(83,15).Value = '15/10/2023 17:20'
$ws.Cells.Item(83,16).Value = 3.49
$ws.Cells.Item(83,17).Value = '20/10/2023 16:59'
$ws.Cells.Item(83,18).Value = 2.25
$ws.Cells.Item(83,19).Value = '15/10/2023 17:20'
$ws.Cells.Item(83,20).Value = 2.45
$ws.Cells.Item(83,21).Value = '20/10/2023 16:59'
$ws.Cells.Item(83,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-akhdoud-al-feiha/EZOH1uVD/'

# Row 84
$ws.Cells.Item(84,6).Value = 'Al Taawon'
$ws.Cells.Item(84,7).Value = 1
$ws.Cells.Item(84,8).Value = 'Al Ittihad'
$ws.Cells.Item(84,9).Value = 1
$ws.Cells.Item(84,10).Value = 4.56
$ws.Cells.Item(84,11).Value = '15/10/2023 12:47'
$ws.Cells.Item(84,12).Value = 4.61
$ws.Cells.Item(84,13).Value = '20/10/2023 16:58'
$ws.Cells.Item(84,14).Value = 4.09
$ws.Cells.Item(84,15).Value = '15/10/2023 12:47'
$ws.Cells.Item(84,16).Value = 4.24
$ws.Cells.Item(84,17).Value = '20/10/2023 16:58'
$ws.Cells.Item(84,18).Value = 1.63
$ws.Cells.Item(84,19).Value = '15/10/2023 12:47'
$ws.Cells.Item(84,20).Value = 1.68
$ws.Cells.Item(84,21).Value = '20/10/2023 16:51'
$ws.Cells.Item(84,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taawon-al-ittihad/44sRvIN0/'

# Row 87
$ws.Cells.Item(87,6).Value = 'Al Nassr'
$ws.Cells.Item(87,7).Value = 2
$ws.Cells.Item(87,8).Value = 'Damac'
$ws.Cells.Item(87,9).Value = 1
$ws.Cells.Item(87,10).Value = 1.17
$ws.Cells.Item(87,11).Value = '14/10/2023 17:02'
$ws.Cells.Item(87,12).Value = 1.16
$ws.Cells.Item(87,13).Value = '21/10/2023 16:59'
$ws.Cells.Item(87,14).Value = 8.119999999999999
$ws.Cells.Item(87,15).Value = '14/10/2023 17:02'
$ws.Cells.Item(87,16).Value = 8.359999999999999
$ws.Cells.Item(87,17).Value = '21/10/2023 16:59'
$ws.Cells.Item(87,18).Value = 13.61
$ws.Cells.Item(87,19).Value = '14/10/2023 17:02'
$ws.Cells.Item(87,20).Value = 14.15
$ws.Cells.Item(87,21).Value = '21/10/2023 16:59'
$ws.Cells.Item(87,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-nassr-damac/QDrNub8f/'

# Row 88
$ws.Cells.Item(88,6).Value = 'Al Hazem'
$ws.Cells.Item(88,7).Value = 4
$ws.Cells.Item(88,8).Value = 'Al Raed'
$ws.Cells.Item(88,9).Value = 3
$ws.Cells.Item(88,10).Value = 2.39
$ws.Cells.Item(88,11).Value = '14/10/2023 17:02'
$ws.Cells.Item(88,12).Value = 3.36
$ws.Cells.Item(88,13).Value = '21/10/2023 16:59'
$ws.Cells.Item(88,14).Value = 3.51
$ws.Cells.Item(88,15).Value = '14/10/2023 17:02'
$ws.Cells.Item(88,16).Value = 3.47
$ws.Cells.Item(88,17).Value = '21/10/2023 16:59'
$ws.Cells.Item(88,18).Value = 2.73
$ws.Cells.Item(88,19).Value = '14/10/2023 17:02'
$ws.Cells.Item(88,20).Value = 2.18
$ws.Cells.Item(88,21).Value = '21/10/2023 16:59'
$ws.Cells.Item(88,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-hazem-rass-al-raed/hQNL0apK/'

# Row 92
$ws.Cells.Item(92,6).Value = 'Al Khaleej'
$ws.Cells.Item(92,7).Value = 1
$ws.Cells.Item(92,8).Value = 'Al Taawon'
$ws.Cells.Item(92,9).Value = 1
$ws.Cells.Item(92,10).Value = 3.53
$ws.Cells.Item(92,11).Value = '24/10/2023 21:01'
$ws.Cells.Item(92,12).Value = 3.5
$ws.Cells.Item(92,13).Value = '26/10/2023 16:55'
$ws.Cells.Item(92,14).Value = 3.77
$ws.Cells.Item(92,15).Value = '24/10/2023 21:01'
$ws.Cells.Item(92,16).Value = 3.7
$ws.Cells.Item(92,17).Value = '26/10/2023 16:56'
$ws.Cells.Item(92,18).Value = 2
$ws.Cells.Item(92,19).Value = '24/10/2023 21:01'
$ws.Cells.Item(92,20).Value = 2.04
$ws.Cells.Item(92,21).Value = '26/10/2023 16:55'
$ws.Cells.Item(92,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-khaleej-al-taawon/Wv18ie76/'

# Row 93
$ws.Cells.Item(93,6).Value = 'Damac'
$ws.Cells.Item(93,7).Value = 2
$ws.Cells.Item(93,8).Value = 'Al Akhdoud'
$ws.Cells.Item(93,9).Value = 0
$ws.Cells.Item(93,10).Value = 1.88
$ws.Cells.Item(93,11).Value = '24/10/2023 21:01'
$ws.Cells.Item(93,12).Value = 1.84
$ws.Cells.Item(93,13).Value = '26/10/2023 16:50'
$ws.Cells.Item(93,14).Value = 3.79
$ws.Cells.Item(93,15).Value = '24/10/2023 21:01'
$ws.Cells.Item(93,16).Value = 3.95
$ws.Cells.Item(93,17).Value = '26/10/2023 17:00'
$ws.Cells.Item(93,18).Value = 3.92
$ws.Cells.Item(93,19).Value = '24/10/2023 21:01'
$ws.Cells.Item(93,20).Value = 3.99
$ws.Cells.Item(93,21).Value = '26/10/2023 16:50'
$ws.Cells.Item(93,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/damac-al-akhdoud/S6fHkZyJ/'

# Row 95
$ws.Cells.Item(95,6).Value = 'Al Taee'
$ws.Cells.Item(95,7).Value = 3
$ws.Cells.Item(95,8).Value = 'Al Riyadh'
$ws.Cells.Item(95,9).Value = 2
$ws.Cells.Item(95,10).Value = 1.84
$ws.Cells.Item(95,11).Value = '24/10/2023 22:01'
$ws.Cells.Item(95,12).Value = 1.75
$ws.Cells.Item(95,13).Value = '27/10/2023 16:59'
$ws.Cells.Item(95,14).Value = 3.78
$ws.Cells.Item(95,15).Value = '24/10/2023 22:01'
$ws.Cells.Item(95,16).Value = 3.85
$ws.Cells.Item(95,17).Value = '27/10/2023 16:59'
$ws.Cells.Item(95,18).Value = 4.13
$ws.Cells.Item(95,19).Value = '24/10/2023 22:01'
$ws.Cells.Item(95,20).Value = 4.65
$ws.Cells.Item(95,21).Value = '27/10/2023 16:58'
$ws.Cells.Item(95,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taee-al-riyadh/pxvV8dxs/'

# Row 96
$ws.Cells.Item(96,6).Value = 'Al Raed'
$ws.Cells.Item(96,7).Value = 1
$ws.Cells.Item(96,8).Value = 'Al Fateh'
$ws.Cells.Item(96,9).Value = 2
$ws.Cells.Item(96,10).Value = 4.14
$ws.Cells.Item(96,11).Value = '24/10/2023 22:01'
$ws.Cells.Item(96,12).Value = 2.93
$ws.Cells.Item(96,13).Value = '27/10/2023 16:59'
$ws.Cells.Item(96,14).Value = 4.26
$ws.Cells.Item(96,15).Value = '24/10/2023 22:01'
$ws.Cells.Item(96,16).Value = 3.98
$ws.Cells.Item(96,17).Value = '27/10/2023 16:59'
$ws.Cells.Item(96,18).Value = 1.75
$ws.Cells.Item(96,19).Value = '24/10/2023 22:01'
$ws.Cells.Item(96,20).Value = 2.21
$ws.Cells.Item(96,21).Value = '27/10/2023 16:59'
$ws.Cells.Item(96,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-raed-al-fateh/lGeDjFMC/'

# ----- New rows appended at the bottom (104:106) -----
# Copy number-format/border/font styling from the last existing data row (103)
# for the styled columns (A: bold+border 'Indice', E: datetime 'data_partida')
# so no duplicate style entries are created in styles.xml.
$ws.Range("A103").Copy() | Out-Null
$ws.Range("A104").PasteSpecial(-4122) | Out-Null
$ws.Range("E103").Copy() | Out-Null
$ws.Range("E104").PasteSpecial(-4122) | Out-Null
$ws.Range("A103").Copy() | Out-Null
$ws.Range("A105").PasteSpecial(-4122) | Out-Null
$ws.Range("E103").Copy() | Out-Null
$ws.Range("E105").PasteSpecial(-4122) | Out-Null
$ws.Range("A103").Copy() | Out-Null
$ws.Range("A106").PasteSpecial(-4122) | Out-Null
$ws.Range("E103").Copy() | Out-Null
$ws.Range("E106").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 104
$ws.Cells.Item(104,1).Value = 103
$ws.Cells.Item(104,2).Value = 'saudi-arabia'
$ws.Cells.Item(104,3).Value = 'saudi-professional-league'
$ws.Cells.Item(104,4).Value = '2023-2024'
$ws.Cells.Item(104,5).Value = 45234.66666666666
$ws.Cells.Item(104,6).Value = 'Abha'
$ws.Cells.Item(104,7).Value = 3
$ws.Cells.Item(104,8).Value = 'Al Akhdoud'
$ws.Cells.Item(104,9).Value = 2
$ws.Cells.Item(104,10).Value = 2.39
$ws.Cells.Item(104,11).Value = '02/11/2023 09:19'
$ws.Cells.Item(104,12).Value = 2.4
$ws.Cells.Item(104,13).Value = '04/11/2023 15:38'
$ws.Cells.Item(104,14).Value = 3.56
$ws.Cells.Item(104,15).Value = '02/11/2023 09:19'
$ws.Cells.Item(104,16).Value = 3.84
$ws.Cells.Item(104,17).Value = '04/11/2023 15:38'
$ws.Cells.Item(104,18).Value = 2.69
$ws.Cells.Item(104,19).Value = '02/11/2023 09:19'
$ws.Cells.Item(104,20).Value = 2.73
$ws.Cells.Item(104,21).Value = '04/11/2023 15:35'
$ws.Cells.Item(104,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/abha-al-akhdoud/0QwZQebP/'

# Row 105
$ws.Cells.Item(105,1).Value = 104
$ws.Cells.Item(105,2).Value = 'saudi-arabia'
$ws.Cells.Item(105,3).Value = 'saudi-professional-league'
$ws.Cells.Item(105,4).Value = '2023-2024'
$ws.Cells.Item(105,5).Value = 45234.79166666666
$ws.Cells.Item(105,6).Value = 'Al Ettifaq'
$ws.Cells.Item(105,7).Value = 0
$ws.Cells.Item(105,8).Value = 'Al Raed'
$ws.Cells.Item(105,9).Value = 0
$ws.Cells.Item(105,10).Value = 1.67
$ws.Cells.Item(105,11).Value = '29/10/2023 19:43'
$ws.Cells.Item(105,12).Value = 2
$ws.Cells.Item(105,13).Value = '04/11/2023 18:55'
$ws.Cells.Item(105,14).Value = 4.1
$ws.Cells.Item(105,15).Value = '29/10/2023 19:43'
$ws.Cells.Item(105,16).Value = 3.42
$ws.Cells.Item(105,17).Value = '04/11/2023 18:55'
$ws.Cells.Item(105,18).Value = 4.82
$ws.Cells.Item(105,19).Value = '29/10/2023 19:43'
$ws.Cells.Item(105,20).Value = 3.96
$ws.Cells.Item(105,21).Value = '04/11/2023 18:55'
$ws.Cells.Item(105,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ettifaq-fc-al-raed/OSxZ7Gim/'

# Row 106
$ws.Cells.Item(106,1).Value = 105
$ws.Cells.Item(106,2).Value = 'saudi-arabia'
$ws.Cells.Item(106,3).Value = 'saudi-professional-league'
$ws.Cells.Item(106,4).Value = '2023-2024'
$ws.Cells.Item(106,5).Value = 45234.79166666666
$ws.Cells.Item(106,6).Value = 'Al Nassr'
$ws.Cells.Item(106,7).Value = 2
$ws.Cells.Item(106,8).Value = 'Al Khaleej'
$ws.Cells.Item(106,9).Value = 0
$ws.Cells.Item(106,10).Value = 1.12
$ws.Cells.Item(106,11).Value = '29/10/2023 19:43'
$ws.Cells.Item(106,12).Value = 1.16
$ws.Cells.Item(106,13).Value = '04/11/2023 18:59'
$ws.Cells.Item(106,14).Value = 10.36
$ws.Cells.Item(106,15).Value = '29/10/2023 19:43'
$ws.Cells.Item(106,16).Value = 8.41
$ws.Cells.Item(106,17).Value = '04/11/2023 18:59'
$ws.Cells.Item(106,18).Value = 17.36
$ws.Cells.Item(106,19).Value = '29/10/2023 19:43'
$ws.Cells.Item(106,20).Value = 13.11
$ws.Cells.Item(106,21).Value = '04/11/2023 18:59'
$ws.Cells.Item(106,22).Value = 'https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-nassr-al-khaleej/EHVo5Ey6/'

